# data : case 1
# Widen columns A and B slightly (stored XML width 15.42578125 -> 16.42578125,
# i.e. +1 character unit) and overwrite the five data rows in A1:B5 with the
# new values from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: A and B widen from a stored XML width of 15.42578125
# to 16.42578125 (one character unit wider). ColumnWidth is expressed in
# "characters", which the file stores with a small fixed padding offset;
# subtract that offset (5/6) so the saved width lands on the target value.
$ws.Columns.Item(1).ColumnWidth = 15.592447916666666
$ws.Columns.Item(2).ColumnWidth = 15.592447916666666

# --- Cell values -----------------------------------------------------------
$ws.Range("A1").Value = -0.00090836526578420677
$ws.Range("B1").Value = -0.00090789938524150904

$ws.Range("A2").Value = -0.026971880409021327
$ws.Range("B2").Value = -0.026970960483861831

$ws.Range("A3").Value = -0.011986400458099251
$ws.Range("B3").Value = -0.012111267906715415

$ws.Range("A4").Value = -0.02030932140230958
$ws.Range("B4").Value = -0.018048881019651273

$ws.Range("A5").Value = -0.030966062255442805
$ws.Range("B5").Value = -0.030966695997759957
